$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "testing responsive" rows with the new test date (10/13/2021 = 44482) ---
#
# A handful of "Start Testing" / "Finish Testing" cells for the Responsive rows
# currently have a plain (bold) style instead of the date-formatted style that
# is already used elsewhere in the sheet (e.g. G17, which uses the non-bold
# "mm-dd-yy" date format). Copy that formatting onto the cells that need it
# before writing the new date values, so the cells end up visually consistent
# with the rest of the date column.
$dateFormatSource = $ws.Range("G17")
$dateFormatSource.Copy()
$ws.Range("F14").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F35").PasteSpecial(-4122)
$ws.Range("F38").PasteSpecial(-4122)
$ws.Range("F48").PasteSpecial(-4122)
$ws.Range("G48").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 14, 28, 35, 38: set "Start Testing" date and bump "% Testing" to 80%.
$ws.Range("F14").Value = 44482
$ws.Range("H14").Value = 0.8

$ws.Range("F28").Value = 44482
$ws.Range("H28").Value = 0.8

$ws.Range("F35").Value = 44482
$ws.Range("H35").Value = 0.8

$ws.Range("F38").Value = 44482
$ws.Range("H38").Value = 0.8

# Row 17 already has the correct date formatting applied, so just set the values.
$ws.Range("F17").Value = 44482
$ws.Range("H17").Value = 0.8

# Rows 44-47: testing is now finished, so set "Finish Testing" date and mark
# "% Testing" complete (100%).
$ws.Range("G44").Value = 44482
$ws.Range("H44").Value = 1

$ws.Range("G45").Value = 44482
$ws.Range("H45").Value = 1

$ws.Range("G46").Value = 44482
$ws.Range("H46").Value = 1

$ws.Range("G47").Value = 44482
$ws.Range("H47").Value = 1

# Row 48: set the new "Start Testing" date and bump "% Testing" to 80%
# (Finish Testing stays blank, only its date formatting was applied above).
$ws.Range("F48").Value = 44482
$ws.Range("H48").Value = 0.8

# --- Update the sheet view ---
# Scroll back to the top of the sheet (clear the previous topLeftCell="A34")
# and move the active selection to I48.
$ws.Range("A1").Select()
$ws.Range("I48").Select()
